# Update column F (dSF) values per repulled data / mean calculation fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -9
$ws.Range("F10").Value = -5
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = -8
$ws.Range("F14").Value = -5
$ws.Range("F16").Value = -7
$ws.Range("F17").Value = -6
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("F22").Value = -3
$ws.Range("F24").Value = -6
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = -4
$ws.Range("F27").Value = -2
$ws.Range("F28").Value = -5
$ws.Range("F29").Value = -3
$ws.Range("F32").Value = -3
$ws.Range("F34").Value = -5
$ws.Range("F35").Value = -10
$ws.Range("F36").Value = -14
$ws.Range("F38").Value = 1
$ws.Range("F39").Value = -6
$ws.Range("F41").Value = 1
$ws.Range("F43").Value = 4
$ws.Range("F44").Value = 6
